$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.75
$ws.Range("S2").Value = 0.25
$ws.Range("B3").Value = 0.3333333333333333
$ws.Range("P3").Value = 0.6666666666666666
$ws.Range("S4").Value = 1
$ws.Range("J6").Value = 0.4
$ws.Range("Q6").Value = 0.2
$ws.Range("S6").Value = 0.4
$ws.Range("J7").Value = 0.2
$ws.Range("Q7").Value = 0.4
$ws.Range("S7").Value = 0.4
$ws.Range("B8").Value = 0.06666666666666667
$ws.Range("F8").Value = 0.06666666666666667
$ws.Range("J8").Value = 0.2
$ws.Range("Q8").Value = 0.3333333333333333
$ws.Range("R8").Value = 0.06666666666666667
$ws.Range("S8").Value = 0.2666666666666667
$ws.Range("Q9").Value = 0.3333333333333333
$ws.Range("R9").Value = 0.5555555555555556
$ws.Range("S9").Value = 0.1111111111111111
$ws.Range("B10").Value = 0.01724137931034483
$ws.Range("D10").Value = 0.01724137931034483
$ws.Range("F10").Value = 0.06896551724137931
$ws.Range("J10").Value = 0.103448275862069
$ws.Range("Q10").Value = 0.3103448275862069
$ws.Range("R10").Value = 0.1206896551724138
$ws.Range("S10").Value = 0.3620689655172414
$ws.Range("G11").Value = 0.2
$ws.Range("J11").Value = 0.1
$ws.Range("K11").Value = 0.3
$ws.Range("L11").Value = 0.4
$ws.Range("G12").Value = 0.75
$ws.Range("J12").Value = 0.25
$ws.Range("I15").Value = 0.1428571428571428
$ws.Range("J15").Value = 0.4285714285714285
$ws.Range("O15").Value = 0.1428571428571428
$ws.Range("S15").Value = 0.2857142857142857
$ws.Range("H16").Value = 0.5
$ws.Range("I16").Value = 0.5
$ws.Range("F17").Value = 0.06896551724137931
$ws.Range("H17").Value = 0.2068965517241379
$ws.Range("I17").Value = 0.1379310344827586
$ws.Range("J17").Value = 0.3103448275862069
$ws.Range("K17").Value = 0.06896551724137931
$ws.Range("M17").Value = 0.03448275862068965
$ws.Range("O17").Value = 0.06896551724137931
$ws.Range("S17").Value = 0.103448275862069
$ws.Range("H18").Value = 0.1538461538461539
$ws.Range("I18").Value = 0.07692307692307693
$ws.Range("J18").Value = 0.3076923076923077
$ws.Range("O18").Value = 0.2307692307692308
$ws.Range("S18").Value = 0.2307692307692308
$ws.Range("F19").Value = 0.0425531914893617
$ws.Range("H19").Value = 0.1276595744680851
$ws.Range("I19").Value = 0.0425531914893617
$ws.Range("J19").Value = 0.574468085106383
$ws.Range("K19").Value = 0.1063829787234043
$ws.Range("S19").Value = 0.1063829787234043
